$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Cells changing type/style (copy from stable donor cells to preserve exact style+shared-string) ---
$ws.Range("C14").Copy($ws.Range("D14"))
$ws.Range("L14").Copy($ws.Range("E14"))
$ws.Range("C14").Copy($ws.Range("G15"))
$ws.Range("L14").Copy($ws.Range("H15"))
$ws.Range("C14").Copy($ws.Range("D20"))
$ws.Range("L14").Copy($ws.Range("E20"))
$ws.Range("I14").Copy($ws.Range("C22"))
$ws.Range("C14").Copy($ws.Range("G26"))
$ws.Range("L14").Copy($ws.Range("H26"))
$ws.Range("C14").Copy($ws.Range("D27"))
$ws.Range("L14").Copy($ws.Range("E27"))
$ws.Range("C14").Copy($ws.Range("G30"))
$ws.Range("L14").Copy($ws.Range("H30"))

# --- Plain value updates (same style/type, value changes only) ---
$ws.Range("G14").Value2 = 1
$ws.Range("F15").Value2 = 1
$ws.Range("D16").Value2 = 3
$ws.Range("E16").Value2 = -33.333333333333
$ws.Range("F16").Value2 = 6
$ws.Range("G16").Value2 = 12
$ws.Range("H16").Value2 = -50
$ws.Range("I16").Value2 = 41
$ws.Range("J16").Value2 = 57
$ws.Range("K16").Value2 = -28.070175438596
$ws.Range("L16").Value2 = 156.25
$ws.Range("M16").Value2 = -4.651162790697
$ws.Range("N16").Value2 = -84.814814814814
$ws.Range("C17").Value2 = 2
$ws.Range("D17").Value2 = 5
$ws.Range("E17").Value2 = -60
$ws.Range("F17").Value2 = 14
$ws.Range("G17").Value2 = 22
$ws.Range("H17").Value2 = -36.363636363636
$ws.Range("I17").Value2 = 65
$ws.Range("J17").Value2 = 69
$ws.Range("K17").Value2 = -5.797101449275
$ws.Range("L17").Value2 = 14.035087719298
$ws.Range("M17").Value2 = 30
$ws.Range("N17").Value2 = -8.450704225352
$ws.Range("C18").Value2 = 1
$ws.Range("D18").Value2 = 1
$ws.Range("F18").Value2 = 11
$ws.Range("H18").Value2 = 37.5
$ws.Range("I18").Value2 = 40
$ws.Range("J18").Value2 = 61
$ws.Range("K18").Value2 = -34.426229508196
$ws.Range("L18").Value2 = 11.111111111111
$ws.Range("M18").Value2 = 122.222222222222
$ws.Range("N18").Value2 = -66.386554621848
$ws.Range("C19").Value2 = 5
$ws.Range("D19").Value2 = 17
$ws.Range("E19").Value2 = -70.588235294117
$ws.Range("F19").Value2 = 44
$ws.Range("G19").Value2 = 64
$ws.Range("H19").Value2 = -31.25
$ws.Range("I19").Value2 = 182
$ws.Range("J19").Value2 = 265
$ws.Range("K19").Value2 = -31.320754716981
$ws.Range("L19").Value2 = 18.181818181818
$ws.Range("M19").Value2 = 152.777777777778
$ws.Range("N19").Value2 = 24.657534246575
$ws.Range("F20").Value2 = 4
$ws.Range("H20").Value2 = -20
$ws.Range("L20").Value2 = -6.25
$ws.Range("M20").Value2 = -34.782608695652
$ws.Range("N20").Value2 = -88.461538461538
$ws.Range("C21").Value2 = 10
$ws.Range("D21").Value2 = 26
$ws.Range("E21").Value2 = -61.538461538461
$ws.Range("F21").Value2 = 80
$ws.Range("H21").Value2 = -28.571428571428
$ws.Range("I21").Value2 = 348
$ws.Range("J21").Value2 = 477
$ws.Range("K21").Value2 = -27.044025157232
$ws.Range("L21").Value2 = 22.968197879858
$ws.Range("M21").Value2 = 66.507177033492
$ws.Range("N21").Value2 = -53.099730458221
$ws.Range("F22").Value2 = 2
$ws.Range("H22").Value2 = 100
$ws.Range("I22").Value2 = 8
$ws.Range("K22").Value2 = -38.461538461538
$ws.Range("L22").Value2 = 60
$ws.Range("M22").Value2 = 166.666666666667
$ws.Range("G23").Value2 = 24
$ws.Range("H23").Value2 = -66.666666666666
$ws.Range("I23").Value2 = 46
$ws.Range("J23").Value2 = 61
$ws.Range("K23").Value2 = -24.590163934426
$ws.Range("L23").Value2 = -25.806451612903
$ws.Range("M23").Value2 = -9.803921568627
$ws.Range("C24").Value2 = 19
$ws.Range("D24").Value2 = 32
$ws.Range("E24").Value2 = -40.625
$ws.Range("F24").Value2 = 82
$ws.Range("G24").Value2 = 173
$ws.Range("H24").Value2 = -52.601156069364
$ws.Range("I24").Value2 = 367
$ws.Range("J24").Value2 = 809
$ws.Range("K24").Value2 = -54.635352286773
$ws.Range("L24").Value2 = 6.686046511627
$ws.Range("M24").Value2 = 60.262008733624
$ws.Range("C25").Value2 = 10
$ws.Range("D25").Value2 = 17
$ws.Range("E25").Value2 = -41.176470588235
$ws.Range("F25").Value2 = 34
$ws.Range("G25").Value2 = 38
$ws.Range("H25").Value2 = -10.526315789473
$ws.Range("I25").Value2 = 138
$ws.Range("J25").Value2 = 160
$ws.Range("K25").Value2 = -13.75
$ws.Range("L25").Value2 = 60.465116279069
$ws.Range("M25").Value2 = 28.971962616822
$ws.Range("F26").Value2 = 1
$ws.Range("C27").Value2 = 2
$ws.Range("G27").Value2 = 5
$ws.Range("H27").Value2 = 40
$ws.Range("I27").Value2 = 20
$ws.Range("K27").Value2 = 33.333333333333
$ws.Range("L27").Value2 = 81.818181818181
$ws.Range("N28").Value2 = -75
$ws.Range("N29").Value2 = -71.428571428571

# --- Header text updates (Volume/Number and date range) - done last to avoid disrupting shared-string indices used above ---
$ws.Range("A8").Value2 = "Volume 30   Number  18"
$ws.Range("C9").Value2 = "Report Covering the Week  5/1/2023  Through  5/7/2023"
